# Results from July 22, 2020 07:43:29 PM America/Los_Angeles TZ run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: New York -- New York -- now has successful data ---
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("B4").Value = 44034

$ws.Range("C4").Value = "'219128"
$ws.Range("D4").Value = "'18803"

$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43

$ws.Range("J4").Value = $true

$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217

$ws.Range("O4").Value = "Success!"

# --- Row 41: Iowa -- updated counts ---
$ws.Range("C41").Value = 40146
$ws.Range("E41").Value = 3289
$ws.Range("G41").Value = 8.19

# --- Row 44: Idaho -- now errored out, data cleared ---
$ws.Range("B44:H44").Clear()
$ws.Range("K44:L44").Clear()

$ws.Range("J44").Value = $false

$ws.Range("O44").Value = "An error occurred. ... TimeoutException('', None, None)"
